$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.051607138645891
$ws.Cells.Item(2, 4).Value = 1.053976467529273
$ws.Cells.Item(2, 5).Value = 1.065328757532329
$ws.Cells.Item(2, 6).Value = 1.072894329116391
$ws.Cells.Item(2, 9).Value = 1.050232964101631
$ws.Cells.Item(2, 10).Value = 1.056634120021556
$ws.Cells.Item(2, 11).Value = 1.056720950572635
$ws.Cells.Item(2, 12).Value = 1.068042292422772
$ws.Cells.Item(2, 13).Value = 1.075587632297378
$ws.Cells.Item(2, 14).Value = 1.058134661669614
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.052613162085364
$ws.Cells.Item(3, 4).Value = 1.054765538461226
$ws.Cells.Item(3, 5).Value = 1.066326987737245
$ws.Cells.Item(3, 6).Value = 1.073973639687318
$ws.Cells.Item(3, 9).Value = 1.050560047362555
$ws.Cells.Item(3, 10).Value = 1.057290019038901
$ws.Cells.Item(3, 11).Value = 1.057323234718844
$ws.Cells.Item(3, 12).Value = 1.068855443175282
$ws.Cells.Item(3, 13).Value = 1.076483127765512
$ws.Cells.Item(3, 14).Value = 1.058791492138796
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.05326438819288
$ws.Cells.Item(4, 4).Value = 1.055276255837048
$ws.Cells.Item(4, 5).Value = 1.066973532967095
$ws.Cells.Item(4, 6).Value = 1.074672790812129
$ws.Cells.Item(4, 9).Value = 1.050770530214602
$ws.Cells.Item(4, 10).Value = 1.057714071214058
$ws.Cells.Item(4, 11).Value = 1.057712430739753
$ws.Cells.Item(4, 12).Value = 1.069381612300978
$ws.Cells.Item(4, 13).Value = 1.077062735524225
$ws.Cells.Item(4, 14).Value = 1.059216146516682
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.053538225949862
$ws.Cells.Item(5, 4).Value = 1.055490993082424
$ws.Cells.Item(5, 5).Value = 1.067245489191698
$ws.Cells.Item(5, 6).Value = 1.074966896492501
$ws.Cells.Item(5, 9).Value = 1.050858738729626
$ws.Cells.Item(5, 10).Value = 1.057892256469196
$ws.Cells.Item(5, 11).Value = 1.057875923211735
$ws.Cells.Item(5, 12).Value = 1.069602814733971
$ws.Cells.Item(5, 13).Value = 1.07730644092917
$ws.Cells.Item(5, 14).Value = 1.05939458481532
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.053584208157725
$ws.Cells.Item(6, 4).Value = 1.055527050243269
$ws.Cells.Item(6, 5).Value = 1.067291160533743
$ws.Cells.Item(6, 6).Value = 1.075016288817735
$ws.Cells.Item(6, 9).Value = 1.050873532988297
$ws.Cells.Item(6, 10).Value = 1.057922169483934
$ws.Cells.Item(6, 11).Value = 1.057903366939948
$ws.Cells.Item(6, 12).Value = 1.069639955634418
$ws.Cells.Item(6, 13).Value = 1.077347362369097
$ws.Cells.Item(6, 14).Value = 1.05942454030997
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.053268046982005
$ws.Cells.Item(7, 4).Value = 1.055279125043265
$ws.Cells.Item(7, 5).Value = 1.066977166276368
$ws.Cells.Item(7, 6).Value = 1.074676719950077
$ws.Cells.Item(7, 9).Value = 1.050771709955707
$ws.Cells.Item(7, 10).Value = 1.057716452473785
$ws.Cells.Item(7, 11).Value = 1.057714615828111
$ws.Cells.Item(7, 12).Value = 1.069384568015737
$ws.Cells.Item(7, 13).Value = 1.077065991777954
$ws.Cells.Item(7, 14).Value = 1.05921853115807
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.051947074589166
$ws.Cells.Item(8, 4).Value = 1.054243109144073
$ws.Cells.Item(8, 5).Value = 1.065665984595062
$ws.Cells.Item(8, 6).Value = 1.073258928415302
$ws.Cells.Item(8, 9).Value = 1.050343743806545
$ws.Cells.Item(8, 10).Value = 1.056855858068508
$ws.Cells.Item(8, 11).Value = 1.056924603095295
$ws.Cells.Item(8, 12).Value = 1.068317098905412
$ws.Cells.Item(8, 13).Value = 1.075890235241991
$ws.Cells.Item(8, 14).Value = 1.058356714610027
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.049621371482729
$ws.Cells.Item(9, 4).Value = 1.052418592704015
$ws.Cells.Item(9, 5).Value = 1.063360319722098
$ws.Cells.Item(9, 6).Value = 1.070766483900575
$ws.Cells.Item(9, 9).Value = 1.049580728021646
$ws.Cells.Item(9, 10).Value = 1.055336657736073
$ws.Cells.Item(9, 11).Value = 1.055528524983474
$ws.Cells.Item(9, 12).Value = 1.066436153889495
$ws.Cells.Item(9, 13).Value = 1.073819669402778
$ws.Cells.Item(9, 14).Value = 1.056835356838864
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.048072278033653
$ws.Cells.Item(10, 4).Value = 1.051203019581938
$ws.Cells.Item(10, 5).Value = 1.061826477750621
$ws.Cells.Item(10, 6).Value = 1.069108844407342
$ws.Cells.Item(10, 9).Value = 1.049066095045526
$ws.Cells.Item(10, 10).Value = 1.054322054173479
$ws.Cells.Item(10, 11).Value = 1.054595164693881
$ws.Cells.Item(10, 12).Value = 1.065182270015744
$ws.Cells.Item(10, 13).Value = 1.072440175176468
$ws.Cells.Item(10, 14).Value = 1.055819312422836
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.047401832656646
$ws.Cells.Item(11, 4).Value = 1.050676855517628
$ws.Cells.Item(11, 5).Value = 1.061163089637022
$ws.Cells.Item(11, 6).Value = 1.068392020938501
$ws.Cells.Item(11, 9).Value = 1.048841844076087
$ws.Cells.Item(11, 10).Value = 1.05388229720118
$ws.Cells.Item(11, 11).Value = 1.054190388894369
$ws.Cells.Item(11, 12).Value = 1.06463934929551
$ws.Cells.Item(11, 13).Value = 1.071843054273945
$ws.Cells.Item(11, 14).Value = 1.055378930945195
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.04715284796037
$ws.Cells.Item(12, 4).Value = 1.050481443586054
$ws.Cells.Item(12, 5).Value = 1.060916794790172
$ws.Cells.Item(12, 6).Value = 1.068125903244618
$ws.Cells.Item(12, 9).Value = 1.048758335440201
$ws.Cells.Item(12, 10).Value = 1.053718888099665
$ws.Cells.Item(12, 11).Value = 1.054039943719763
$ws.Cells.Item(12, 12).Value = 1.064437687671956
$ws.Cells.Item(12, 13).Value = 1.071621288811012
$ws.Cells.Item(12, 14).Value = 1.055215289784012
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.047206253822722
$ws.Cells.Item(13, 4).Value = 1.050523358792616
$ws.Cells.Item(13, 5).Value = 1.060969620566497
$ws.Cells.Item(13, 6).Value = 1.068182979939839
$ws.Cells.Item(13, 9).Value = 1.048776257908768
$ws.Cells.Item(13, 10).Value = 1.053753942787099
$ws.Cells.Item(13, 11).Value = 1.054072218941036
$ws.Cells.Item(13, 12).Value = 1.064480944618571
$ws.Cells.Item(13, 13).Value = 1.071668856808961
$ws.Cells.Item(13, 14).Value = 1.055250394253123
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.047381250512344
$ws.Cells.Item(14, 4).Value = 1.050660702114147
$ws.Cells.Item(14, 5).Value = 1.061142728437673
$ws.Cells.Item(14, 6).Value = 1.0683700206583
$ws.Cells.Item(14, 9).Value = 1.048834945546404
$ws.Cells.Item(14, 10).Value = 1.053868791052911
$ws.Cells.Item(14, 11).Value = 1.054177954952584
$ws.Cells.Item(14, 12).Value = 1.064622679812985
$ws.Cells.Item(14, 13).Value = 1.071824722420027
$ws.Cells.Item(14, 14).Value = 1.055365405616646
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.047489078188646
$ws.Cells.Item(15, 4).Value = 1.050745327709067
$ws.Cells.Item(15, 5).Value = 1.06124940143419
$ws.Cells.Item(15, 6).Value = 1.068485281508134
$ws.Cells.Item(15, 9).Value = 1.048871076872525
$ws.Cells.Item(15, 10).Value = 1.05393954441906
$ws.Cells.Item(15, 11).Value = 1.054243090032695
$ws.Cells.Item(15, 12).Value = 1.064710007994447
$ws.Cells.Item(15, 13).Value = 1.071920760596554
$ws.Cells.Item(15, 14).Value = 1.05543625946069
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.04811678008821
$ws.Cells.Item(16, 4).Value = 1.051237943361569
$ws.Cells.Item(16, 5).Value = 1.061870521019288
$ws.Cells.Item(16, 6).Value = 1.069156437574466
$ws.Cells.Item(16, 9).Value = 1.049080948123962
$ws.Cells.Item(16, 10).Value = 1.054351230431102
$ws.Cells.Item(16, 11).Value = 1.054622015206432
$ws.Cells.Item(16, 12).Value = 1.065218302316852
$ws.Cells.Item(16, 13).Value = 1.072479808565628
$ws.Cells.Item(16, 14).Value = 1.05584853011409
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.048510607298407
$ws.Cells.Item(17, 4).Value = 1.051546998748777
$ws.Cells.Item(17, 5).Value = 1.06226034093231
$ws.Cells.Item(17, 6).Value = 1.069577689446682
$ws.Cells.Item(17, 9).Value = 1.049212217072609
$ws.Cells.Item(17, 10).Value = 1.054609356141067
$ws.Cells.Item(17, 11).Value = 1.054859538023408
$ws.Cells.Item(17, 12).Value = 1.06553714736465
$ws.Cells.Item(17, 13).Value = 1.072830540875125
$ws.Cells.Item(17, 14).Value = 1.056107022392171
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.048740351135247
$ws.Cells.Item(18, 4).Value = 1.051727283533534
$ws.Cells.Item(18, 5).Value = 1.06248779109018
$ws.Cells.Item(18, 6).Value = 1.069823489722404
$ws.Cells.Item(18, 9).Value = 1.049288647825783
$ws.Cells.Item(18, 10).Value = 1.054759875283295
$ws.Cells.Item(18, 11).Value = 1.054998020767226
$ws.Cells.Item(18, 12).Value = 1.065723126118369
$ws.Cells.Item(18, 13).Value = 1.073035137255379
$ws.Cells.Item(18, 14).Value = 1.056257755288846
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.048818693115243
$ws.Cells.Item(19, 4).Value = 1.051788759018518
$ws.Cells.Item(19, 5).Value = 1.062565358387601
$ws.Cells.Item(19, 6).Value = 1.069907316718083
$ws.Cells.Item(19, 9).Value = 1.049314685628305
$ws.Cells.Item(19, 10).Value = 1.054811191418945
$ws.Cells.Item(19, 11).Value = 1.055045229563895
$ws.Cells.Item(19, 12).Value = 1.065786540365527
$ws.Cells.Item(19, 13).Value = 1.073104902740923
$ws.Cells.Item(19, 14).Value = 1.056309144299295
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.048468350149166
$ws.Cells.Item(20, 4).Value = 1.051513838162902
$ws.Cells.Item(20, 5).Value = 1.062218509204389
$ws.Cells.Item(20, 6).Value = 1.069532483672549
$ws.Cells.Item(20, 9).Value = 1.049198147243296
$ws.Cells.Item(20, 10).Value = 1.0545816659559
$ws.Cells.Item(20, 11).Value = 1.054834060306699
$ws.Cells.Item(20, 12).Value = 1.065502938116044
$ws.Cells.Item(20, 13).Value = 1.072792908529633
$ws.Cells.Item(20, 14).Value = 1.056079292883765
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.047329717007015
$ws.Cells.Item(21, 4).Value = 1.050620257127558
$ws.Cells.Item(21, 5).Value = 1.06109174924931
$ws.Cells.Item(21, 6).Value = 1.06831493789256
$ws.Cells.Item(21, 9).Value = 1.04881766934594
$ws.Cells.Item(21, 10).Value = 1.053834972858932
$ws.Cells.Item(21, 11).Value = 1.054146820908016
$ws.Cells.Item(21, 12).Value = 1.064580942236331
$ws.Cells.Item(21, 13).Value = 1.071778823000628
$ws.Cells.Item(21, 14).Value = 1.055331539396952
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.04661409423521
$ws.Cells.Item(22, 4).Value = 1.050058594417527
$ws.Cells.Item(22, 5).Value = 1.060383987934478
$ws.Cells.Item(22, 6).Value = 1.067550242974868
$ws.Cells.Item(22, 9).Value = 1.048577222254689
$ws.Cells.Item(22, 10).Value = 1.053365128430904
$ws.Cells.Item(22, 11).Value = 1.053714185933626
$ws.Cells.Item(22, 12).Value = 1.064001266162554
$ws.Cells.Item(22, 13).Value = 1.071141411520291
$ws.Cells.Item(22, 14).Value = 1.054861027735945
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.046993432714015
$ws.Cells.Item(23, 4).Value = 1.0503563263635
$ws.Cells.Item(23, 5).Value = 1.06075912121785
$ws.Cells.Item(23, 6).Value = 1.067955543928217
$ws.Cells.Item(23, 9).Value = 1.048704803899633
$ws.Cells.Item(23, 10).Value = 1.053614236692985
$ws.Cells.Item(23, 11).Value = 1.053943584990178
$ws.Cells.Item(23, 12).Value = 1.064308561447177
$ws.Cells.Item(23, 13).Value = 1.071479297799444
$ws.Cells.Item(23, 14).Value = 1.055110489760331
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.048487444239814
$ws.Cells.Item(24, 4).Value = 1.05152882194866
$ws.Cells.Item(24, 5).Value = 1.062237410930467
$ws.Cells.Item(24, 6).Value = 1.069552909934219
$ws.Cells.Item(24, 9).Value = 1.049204505214352
$ws.Cells.Item(24, 10).Value = 1.054594178085987
$ws.Cells.Item(24, 11).Value = 1.054845572776491
$ws.Cells.Item(24, 12).Value = 1.065518395797271
$ws.Cells.Item(24, 13).Value = 1.072809912904363
$ws.Cells.Item(24, 14).Value = 1.056091822782512
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.05022238039784
$ws.Cells.Item(25, 4).Value = 1.052890141510291
$ws.Cells.Item(25, 5).Value = 1.063955816113539
$ws.Cells.Item(25, 6).Value = 1.071410139698424
$ws.Cells.Item(25, 9).Value = 1.049779036952518
$ws.Cells.Item(25, 10).Value = 1.055729727306706
$ws.Cells.Item(25, 11).Value = 1.055889912493253
$ws.Cells.Item(25, 12).Value = 1.066922411540279
$ws.Cells.Item(25, 13).Value = 1.074354806939047
$ws.Cells.Item(25, 14).Value = 1.057228984613373
